$d = $word.ActiveDocument

# 1) Fix the "Deumeurant" typo and merge it with the following "[ADDR]" run
#    into a single run (matches the template's other single-run fields).
$d.Content.Find.Execute("Deumeurant au : [ADDR]", $true, $false, $false, $false, $false, $true, 1, $false, "Demeurant au : [ADDR]", 2)

# 2) Add the missing "ListLabel 22" character style (same shape as the
#    existing ListLabel20 / ListLabel21 styles: Symbol complex-script font,
#    24 half-point size, quick style).
$listLabel22 = $d.Styles.Add("ListLabel 22", 2)
$listLabel22.Font.NameBi = "Symbol"
$listLabel22.Font.Size = 12
$listLabel22.Font.SizeBi = 12
$listLabel22.QuickStyle = $true
